# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 3085
    4  = 193
    5  = 164
    7  = 1705
    9  = 91
    12 = 1401
    13 = 14
    14 = 535
    23 = 3273
    25 = 160
    26 = 338
    27 = 13
    29 = 148
    30 = 105
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
